# Apply the "manual testing spreadsheet completed" edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark options 1-3 (rows 14-16) as complete
$ws.Range("D14").Value = "complete"
$ws.Range("D15").Value = "complete"
$ws.Range("D16").Value = "complete"

# Row 17 (option 4 - call med record method): mark complete, add date
$ws.Range("D17").Value = "complete"
$ws.Range("E17").Value = 43712
$ws.Range("E4").Copy()
$ws.Range("E17").PasteSpecial(-4122)

# Row 18 (ARGV / add medication): mark not implemented
$ws.Range("D18").Value = "not implemented"

# Row 22 (CSV check for view results method): mark partial, add date, add comment
$ws.Range("D22").Value = "partial"
$ws.Range("E22").Value = 43713
$ws.Range("E4").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F17").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F22").Value = "if csv does not exist, rescue error message will display after user types in a name"

# Row 17 comment update (done after F22 so shared-string ordering matches upstream edit order)
$ws.Range("F17").Value = "grab patient name & score from previous entries & use to create instance object of a patient class if possible"

# Restore the saved view/selection state
$ws.Range("D21").Select()
$excel.ActiveWindow.ScrollRow = 2
